$wb = $excel.ActiveWorkbook

# --- 1. Refresh the "panel query" timestamps recorded in the data sheet ---
# (the panel was re-fetched, giving a later time_taken value for every gene row)
$dataSheet = $wb.Worksheets.Item("data")

$timestamps = @(
    "2021-10-05 14:35:50.179543",
    "2021-10-05 14:35:50.179551",
    "2021-10-05 14:35:50.179554",
    "2021-10-05 14:35:50.179557",
    "2021-10-05 14:35:50.179560",
    "2021-10-05 14:35:50.179563",
    "2021-10-05 14:35:50.179565",
    "2021-10-05 14:35:50.179568",
    "2021-10-05 14:35:50.179571",
    "2021-10-05 14:35:50.179573",
    "2021-10-05 14:35:50.179576",
    "2021-10-05 14:35:50.179578",
    "2021-10-05 14:35:50.179581",
    "2021-10-05 14:35:50.179583",
    "2021-10-05 14:35:50.179586",
    "2021-10-05 14:35:50.179588",
    "2021-10-05 14:35:50.179591",
    "2021-10-05 14:35:50.179594",
    "2021-10-05 14:35:50.179597",
    "2021-10-05 14:35:50.179599",
    "2021-10-05 14:35:50.179602",
    "2021-10-05 14:35:50.179604",
    "2021-10-05 14:35:50.179607",
    "2021-10-05 14:35:50.179609",
    "2021-10-05 14:35:50.179612",
    "2021-10-05 14:35:50.179615",
    "2021-10-05 14:35:50.179617",
    "2021-10-05 14:35:50.179620",
    "2021-10-05 14:35:50.179622",
    "2021-10-05 14:35:50.179625",
    "2021-10-05 14:35:50.179628",
    "2021-10-05 14:35:50.179630",
    "2021-10-05 14:35:50.179633",
    "2021-10-05 14:35:50.179636",
    "2021-10-05 14:35:50.179638",
    "2021-10-05 14:35:50.179641",
    "2021-10-05 14:35:50.179643",
    "2021-10-05 14:35:50.179646",
    "2021-10-05 14:35:50.179648",
    "2021-10-05 14:35:50.179651",
    "2021-10-05 14:35:50.179654",
    "2021-10-05 14:35:50.179656",
    "2021-10-05 14:35:50.179659",
    "2021-10-05 14:35:50.179662",
    "2021-10-05 14:35:50.179664",
    "2021-10-05 14:35:50.179667",
    "2021-10-05 14:35:50.179669",
    "2021-10-05 14:35:50.179672",
    "2021-10-05 14:35:50.179674",
    "2021-10-05 14:35:50.179677",
    "2021-10-05 14:35:50.179679",
    "2021-10-05 14:35:50.179682",
    "2021-10-05 14:35:50.179685",
    "2021-10-05 14:35:50.179687",
    "2021-10-05 14:35:50.179690",
    "2021-10-05 14:35:50.179693",
    "2021-10-05 14:35:50.179695"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $timestamps[$i]
}

# --- 2. Add a new "metadata" worksheet, positioned right after "data" ---
# It records the details of the PanelApp query that produced this export.
$metaSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$metaSheet.Name = "metadata"

$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Systemic Autoinflammatory Disease_Periodic Fever"
$metaSheet.Range("C2").Value = 238
$metaSheet.Range("D2").Value = "'0.121"
$metaSheet.Range("E2").Value = "2021-09-10T06:11:16.816747Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:35:50.175841"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/238/?format=json"

# Match the formatting used for the header row / index column on "data"
# (bold + bordered + centered header cells, plain index cells) by copying
# the already-established cell formats across instead of inventing new ones.
$dataSheet.Range("B1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats

$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

# D2 only needed the quote-prefix trick to force text storage of "0.121";
# strip the implicit quote-prefix formatting back off so the cell carries no
# explicit style, same as the other plain data cells.
$metaSheet.Range("C2").Copy()
$metaSheet.Range("D2").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false
